# Auto-generated Excel COM-interop script implementing the diff
$wb = $excel.ActiveWorkbook

# ---- Simple numeric ("想去人数") updates: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 267
$ws.Range("F7").Value = 2174
$ws.Range("F8").Value = 51
$ws.Range("F9").Value = 55
$ws.Range("F10").Value = 1633
$ws.Range("F11").Value = 1633
$ws.Range("F12").Value = 1362
$ws.Range("F13").Value = 63
$ws.Range("F14").Value = 1410
$ws.Range("F17").Value = 579
$ws.Range("F18").Value = 161
$ws.Range("F19").Value = 12
$ws.Range("F20").Value = 7243
$ws.Range("F21").Value = 7953
$ws.Range("F24").Value = 199
$ws.Range("F35").Value = 1438
$ws.Range("F36").Value = 207
$ws.Range("F37").Value = 227
$ws.Range("F40").Value = 12
$ws.Range("F41").Value = 727
$ws.Range("F44").Value = 343
$ws.Range("F45").Value = 245
$ws.Range("F46").Value = 194
$ws.Range("F47").Value = 86
$ws.Range("F48").Value = 178

# ---- Simple numeric ("想去人数") updates: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 19
$ws.Range("F5").Value = 59
$ws.Range("F18").Value = 298

# ---- Simple numeric ("想去人数") updates: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 184
$ws.Range("F3").Value = 2622
$ws.Range("F6").Value = 13

# ---- 全部类型: simple numeric updates unrelated to the new row ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 267
$ws.Range("F3").Value = 184
$ws.Range("F5").Value = 19
$ws.Range("F10").Value = 2174
$ws.Range("F11").Value = 51
$ws.Range("F12").Value = 55
$ws.Range("F13").Value = 1633
$ws.Range("F14").Value = 1633
$ws.Range("F15").Value = 63
$ws.Range("F16").Value = 1410
$ws.Range("F18").Value = 579
$ws.Range("F20").Value = 161
$ws.Range("F21").Value = 59
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 7243
$ws.Range("F25").Value = 7953
$ws.Range("F39").Value = 727
$ws.Range("F44").Value = 343
$ws.Range("F45").Value = 245
$ws.Range("F46").Value = 194
$ws.Range("F47").Value = 178
$ws.Range("F49").Value = 298

# ---- 全部类型: new row (2024-10-01 爱川こずえ) inserted at row 27, shifting rows
# 27-37 down to 28-38; the event previously at row 38 ("黑白键上的音乐地图")
# falls out of this sheet's list (rows 39+ are unaffected). Column A (the
# running index) is not part of the shift, so it is left untouched.
# row 27
$ws.Range("B27").Value = '2024-10-01'
$ws.Range("C27").Value = '北京·超人气舞见 爱川こずえ 签售会'
$ws.Range("D27").Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws.Range("E27").Value = '2024.10.01 11:00-10.01 17:00'
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 158
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=92039'
$ws.Range("I27").Value = '//i1.hdslb.com/bfs/openplatform/202409/18ftTPgv1725605173522.png'

# row 28
$ws.Range("B28").Value = '2024-10-02'
$ws.Range("C28").Value = '北京·人气声优 内田秀 专场活动'
$ws.Range("D28").Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws.Range("E28").Value = '2024.10.02 13:55-10.02 17:10'
$ws.Range("F28").Value = 90
$ws.Range("G28").Value = 458
$ws.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=91678'
$ws.Range("I28").Value = '//i0.hdslb.com/bfs/openplatform/202409/0aUkHD511725260741169.png'

# row 29
$ws.Range("B29").Value = '2024-10-02'
$ws.Range("C29").Value = '北京·第19届IJOY漫展【217专场见面会】'
$ws.Range("D29").Value = '天辰东路7号 北京国家会议中心'
$ws.Range("E29").Value = '2024.10.02 12:25-10.02 16:30'
$ws.Range("F29").Value = 16
$ws.Range("G29").Value = 168
$ws.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=91561'
$ws.Range("I29").Value = '//i0.hdslb.com/bfs/openplatform/202408/cAghXlck1724928163645.jpeg'

# row 30
$ws.Range("B30").Value = '2024-10-02'
$ws.Range("C30").Value = '北京·第19届IJOY漫展【银发娘专场见面会】'
$ws.Range("D30").Value = '天辰东路7号 北京国家会议中心'
$ws.Range("E30").Value = '2024.10.02 12:25-10.02 16:30'
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 168
$ws.Range("H30").Value = 'https://show.bilibili.com/platform/detail.html?id=91563'
$ws.Range("I30").Value = '//i0.hdslb.com/bfs/openplatform/202408/9Td79pPE1724928479521.jpeg'

# row 31
$ws.Range("B31").Value = '2024-10-04'
$ws.Range("C31").Value = '北京·第五人格only同人展'
$ws.Range("D31").Value = '北花园路1号 超级蜂巢'
$ws.Range("E31").Value = '2024.10.04 10:00-10.04 17:00'
$ws.Range("F31").Value = 1438
$ws.Range("G31").Value = 68
$ws.Range("H31").Value = 'https://show.bilibili.com/platform/detail.html?id=89309'
$ws.Range("I31").Value = '//i0.hdslb.com/bfs/openplatform/202407/4XsICpa71721046044404.jpeg'

# row 32
$ws.Range("B32").Value = '2024-10-05'
$ws.Range("C32").Value = '北京·咒术回战同人Only2.0'
$ws.Range("D32").Value = '安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园'
$ws.Range("E32").Value = '2024.10.05 09:30-10.05 17:00'
$ws.Range("F32").Value = 207
$ws.Range("G32").Value = 65
$ws.Range("H32").Value = 'https://show.bilibili.com/platform/detail.html?id=91628'
$ws.Range("I32").Value = '//i0.hdslb.com/bfs/openplatform/202408/IsJo7aU61724405528082.jpeg'

# row 33
$ws.Range("B33").Value = '2024-10-05'
$ws.Range("C33").Value = '北京·马娘ONLY2'
$ws.Range("D33").Value = '永外高庄138号 北京大红门国际会展中心'
$ws.Range("E33").Value = '2024.10.05 10:00-10.05 17:00'
$ws.Range("F33").Value = 227
$ws.Range("G33").Value = 75
$ws.Range("H33").Value = 'https://show.bilibili.com/platform/detail.html?id=89334'
$ws.Range("I33").Value = '//i1.hdslb.com/bfs/openplatform/202408/30C9r9Qz1724639124911.png'

# row 34
$ws.Range("B34").Value = '2024-10-06'
$ws.Range("C34").Value = '北京·Hi Fun 全忍界秋季运动会 火影同人ONLY x 北投购物公园潮街 '
$ws.Range("D34").Value = '安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园'
$ws.Range("E34").Value = '2024.10.06 11:00-10.06 18:00'
$ws.Range("F34").Value = 8
$ws.Range("G34").Value = 68
$ws.Range("H34").Value = 'https://show.bilibili.com/platform/detail.html?id=91930'
$ws.Range("I34").Value = '//i2.hdslb.com/bfs/openplatform/202409/f7nTqmEI1725439502652.jpeg'

# row 35
$ws.Range("B35").Value = '2024-10-06'
$ws.Range("C35").Value = '北京·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会'
$ws.Range("D35").Value = '神路街39号地下一层 DDC 黄昏黎明俱乐部'
$ws.Range("E35").Value = '2024.10.06 19:30-10.06 21:00'
$ws.Range("F35").Value = 5
$ws.Range("G35").Value = 220
$ws.Range("H35").Value = 'https://show.bilibili.com/platform/detail.html?id=91353'
$ws.Range("I35").Value = '//i1.hdslb.com/bfs/openplatform/202408/gwqdCudt1724743063464.jpeg'

# row 36
$ws.Range("B36").Value = '2024-10-06'
$ws.Range("C36").Value = '北京·第七届璃樱动漫嘉年华'
$ws.Range("D36").Value = '永外高庄138号 北京大红门国际会展中心'
$ws.Range("E36").Value = '2024.10.06 10:00-10.06 17:00'
$ws.Range("F36").Value = 291
$ws.Range("G36").Value = 60
$ws.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=90346'
$ws.Range("I36").Value = '//i0.hdslb.com/bfs/openplatform/202408/E10ytIKK1722569456045.jpeg'

# row 37
$ws.Range("B37").Value = '2024-10-07'
$ws.Range("C37").Value = '北京·秋日物语-运动番同人Only'
$ws.Range("D37").Value = '酒仙桥北路2号院798艺术区706后街1号 北京格瑞斯艺术酒店'
$ws.Range("E37").Value = '2024.10.07 10:00-10.07 17:00'
$ws.Range("F37").Value = 12
$ws.Range("G37").Value = 69
$ws.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=92076'
$ws.Range("I37").Value = '//i1.hdslb.com/bfs/openplatform/202409/81nlLDGH1725604082584.png'

# row 38
$ws.Range("B38").Value = '2024-10-07'
$ws.Range("C38").Value = '北京·集结 - 超级世代！ACGN BAND LIVE S8'
$ws.Range("D38").Value = '隆福寺街95号钱粮胡同38号15号楼B1层 東市OMNICLUB'
$ws.Range("E38").Value = '2024.10.07 12:30-10.07 20:30'
$ws.Range("F38").Value = 26
$ws.Range("G38").Value = 88
$ws.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=91317'
$ws.Range("I38").Value = '//i2.hdslb.com/bfs/openplatform/202408/B0Cuvd5v1724740500595.jpeg'

